$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: update title and link
$ws.Range("D9").Value = "ChatGPT 시리즈 – ②’인간 피드백형 강화학습(RLwHF)’의 장점"
$ws.Range("E9").Value = "https://pdsi.pabii.com/chatgpt-series-2/#utm_source=rss&utm_medium=rss&utm_campaign=chatgpt-series-2"

# Row 23: update title and link
$ws.Range("D23").Value = "[공개] 대용량 데이터셋 다운로드 받는 코드(인터넷 끊길 때 이용하면 좋음)"
$ws.Range("E23").Value = "https://theonly1.tistory.com/3125"

# Row 42: update title and link
$ws.Range("D42").Value = "[임베디드]Zynq 7000 TRM(UG585) - 7. Interrupt - 번역"
$ws.Range("E42").Value = "https://kjk92.tistory.com/103"
